$d = $word.ActiveDocument
$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Paragraph 1 (Heading1): "ContextFreeSQL" -> wrap run with proofErr spellStart/spellEnd ---
$p1 = $d.Paragraphs(1)
$p1.Range.InsertXML(@"
<w:p $w w14:paraId="398F5278" w14:textId="56764519" w:rsidR="005E73B4" w:rsidRDefault="001E7502" w:rsidP="001E7502" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>ContextFreeSQL</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
"@)

# --- Paragraph 2: "Smallie: title of server.dbname and timestamp ( as usual, will be in comments or printed out)" ---
# Text content is unchanged; split into multiple runs with spelling/grammar proofErr marks.
$p2 = $d.Paragraphs(2)
$p2.Range.InsertXML(@"
<w:p $w w14:paraId="7DA8A916" w14:textId="337DEA3E" w:rsidR="000A39FE" w:rsidRDefault="000A39FE" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Smallie: title of </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>server.dbname</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> and timestamp </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>( as</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> usual, will be in comments or printed out)</w:t></w:r></w:p>
"@)

# --- Paragraph 3: "Massive testing..." -> "Now: Testing seems ok. Now test the exe, on its own with all command lines, begin with --help" ---
$p3 = $d.Paragraphs(3)
$p3.Range.InsertXML(@"
<w:p $w w14:paraId="79C017E8" w14:textId="036F0EFD" w:rsidR="00D265AC" w:rsidRDefault="00D265AC" w:rsidP="00A346EB" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr></w:pPr><w:r><w:t>Now: Testing seems ok. Now test the exe, on its own with all command lines, begin with --help</w:t></w:r></w:p>
"@)

# --- Paragraph 4: "Now: begin. Set up the website..." -> "Copy to hostinger, test downloads on Windows and Linux" (ilvl 1 -> 0) ---
$p4 = $d.Paragraphs(4)
$p4.Range.InsertXML(@"
<w:p $w w14:paraId="3059CE48" w14:textId="50B7AD8C" w:rsidR="00575A61" w:rsidRDefault="00575A61" w:rsidP="00575A61" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Copy to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>hostinger</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, test downloads on Windows and Linux</w:t></w:r></w:p>
"@)

# --- Paragraph 5: "Is the "-help" good enough?..." -> bold "Now" + ": the movie. How do we do with avatar on camtasia" (ilvl 1 -> 0) ---
$p5 = $d.Paragraphs(5)
$p5.Range.InsertXML(@"
<w:p $w w14:paraId="0B188F59" w14:textId="264E2D37" w:rsidR="00575A61" w:rsidRDefault="00575A61" w:rsidP="00575A61" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Now</w:t></w:r><w:r><w:t xml:space="preserve">: the movie. How do we do with avatar on </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>camtasia</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
"@)

# --- Paragraph 6: "Get claude to generate full doc, for the website" -> removed entirely ---
$p6 = $d.Paragraphs(6)
$p6.Range.Delete()
